# Adds other test logins to the UserList sheet.
# Rows 2-4 previously held single-letter placeholder logins ("a","b","c").
# Replace them with three real test accounts (username, first name, last name).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("UserList")

$ws.Range("B2").Value = "Jack"
$ws.Range("C2").Value = "Edwards"

$ws.Range("B3").Value = "Matthew"
$ws.Range("C3").Value = "Columbo"

$ws.Range("B4").Value = "Angel"
$ws.Range("C4").Value = "Rivera"

$ws.Range("A2").Value = "jedwards"
$ws.Range("A3").Value = "mcolumbo"
$ws.Range("A4").Value = "arivera"

$ws.Range("I16").Select()
